$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2 through 44
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 68   # AD
    $ws.Cells.Item($r, 31).Value = 94   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
